# Update the "想去人数" (number of people wanting to go) column (F) values
# for both the "展览" (exhibitions) sheet and the "全部类型" (all types) sheet,
# reflecting refreshed counts from the latest data pull.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibitions) ---
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F4").Value = 680
$wsExhibit.Range("F5").Value = 577
$wsExhibit.Range("F6").Value = 323
$wsExhibit.Range("F7").Value = 2825
$wsExhibit.Range("F9").Value = 8039
$wsExhibit.Range("F11").Value = 475
$wsExhibit.Range("F12").Value = 48
$wsExhibit.Range("F13").Value = 399

# --- Sheet "全部类型" (All Types) ---
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F4").Value = 680
$wsAll.Range("F5").Value = 577
$wsAll.Range("F6").Value = 323
$wsAll.Range("F9").Value = 2825
$wsAll.Range("F11").Value = 8039
$wsAll.Range("F13").Value = 475
$wsAll.Range("F14").Value = 48
$wsAll.Range("F17").Value = 399

$wb.Save()
